# Generate Report for Handback
# Updates the "Generate"/"Handoff"/"Handback" timestamp cells to new values,
# matching the regenerated report timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for first file row
$wsOverview.Range("G2").Value = "2016-08-27 17:03:33"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime (row 2)
$wsZhCn.Range("H2").Value = "2016-08-27 17:03:29"
$wsZhCn.Range("K2").Value = "2016-08-27 17:03:46"

# de-de sheet: Correspond Handoff Datetime (row 2, shared with Overview!G2 value)
#              Correspond Handback DateTime (row 2)
$wsDeDe.Range("H2").Value = "2016-08-27 17:03:33"
$wsDeDe.Range("K2").Value = "2016-08-27 17:03:53"
